{"js": "// Applies the Question2_Soln.docx edit:\n//  1. \"int(\" -> \"int (\" (first occurrence)\n//  2. \", ie. \" -> \", i.e. \" (congestion-window paragraph)\n//  3. \"...the window, ie. if the window is 5 ...\" -> \"...the window, i.e. if the window is 5 ...\"\n//  4. \"1 / int(2.5)\" -> \"1 / int (2.5)\"\n//  5. \"Window size at time 21 is 2.75\" -> \"Window size at time 21 is 3\"\n//  6. Empty \"Ans:\" paragraph (time-23 question) gets the full answer appended\n//  7. \"Window size at time31 is 3.5.\" -> \"Window size at time 31 is 3.99\" (and\n//     the _GoBack bookmark moves from the trailing empty paragraph into this run)\n\nconst body = context.document.body;\n\n// ---- 1. \"int(\" -> \"int (\" in the first occurrence (\"fractionally by (1 / int(current window size))\") ----\n{\n  const results = body.search(\"int(current window size\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"int (current window size\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// ---- 2. \", ie. \" -> \", i.e. \" in \"...within the window, ie. even with a window...\" ----\n{\n  const results = body.search(\", ie. even with a window\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\", i.e. even with a window\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// ---- 3. \"...increment the window, ie. if the window is 5 ...\" -> \"...i.e. if the window is 5 ...\" ----\n{\n  const results = body.search(\"the window, ie. if the window is 5\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"the window, i.e. if the window is 5\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// ---- 4. \"1 / int(2.5)\" -> \"1 / int (2.5)\" ----\n{\n  const results = body.search(\"1 / int(2.5)\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"1 / int (2.5)\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// ---- 5. \"Window size at time 21 is 2.75\" -> \"Window size at time 21 is 3\" ----\n{\n  const results = body.search(\"Window size at time 21 is 2.75\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"Window size at time 21 is 3\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// ---- 6. Fill in the empty \"Ans:\" paragraph (the one right after the \"time 23\" question) ----\n{\n  const paras = body.paragraphs;\n  paras.load(\"items/text\");\n  await context.sync();\n\n  let target = null;\n  for (let i = 0; i < paras.items.length; i++) {\n    if (paras.items[i].text === \"Ans:\") {\n      target = paras.items[i];\n      break;\n    }\n  }\n\n  if (target) {\n    // The paragraph (and its trailing mark) stops being entirely bold once\n    // non-bold content follows \"Ans: \" - set the whole paragraph/mark to\n    // non-bold first, then re-bold just the existing \"Ans:\" text.\n    target.font.bold = false;\n    await context.sync();\n\n    const ansRuns = target.search(\"Ans:\", { matchCase: true });\n    ansRuns.load(\"items\");\n    await context.sync();\n    if (ansRuns.items.length > 0) {\n      ansRuns.items[0].font.bold = true;\n      await context.sync();\n    }\n\n    const spaceRange = target.insertText(\" \", Word.InsertLocation.end);\n    spaceRange.font.bold = true;\n    await context.sync();\n\n    const answerText =\n      \"At time 23, the window size is 3.33 and the number of packets in the network is 2. \" +\n      \"Therefore, one more packet can be pushed into the network as that can be fully accommodated in the congestion window. \" +\n      \"That\\u2019s the reason why a packet is pushed out at time 23 even though no ack is received at that time.\";\n    const answerRange = target.insertText(answerText, Word.InsertLocation.end);\n    answerRange.font.bold = false;\n    await context.sync();\n  }\n}\n\n// ---- 7. \"Window size at time31 is 3.5.\" -> \"Window size at time 31 is 3.99\" + move _GoBack bookmark ----\n{\n  const results = body.search(\"Window size at time31 is 3.5.\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"Window size at time 31 is 3.99\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n\n  // Relocate the _GoBack bookmark from the trailing empty paragraph to just\n  // before \"31 is 3.99\".\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n\n  const target = body.search(\"31 is 3.99\", { matchCase: true });\n  target.load(\"items\");\n  await context.sync();\n  if (target.items.length > 0) {\n    const collapsed = target.items[0].getRange(Word.RangeLocation.start);\n    collapsed.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "# Applies the Question2_Soln.docx edit:\n#  1. \"int(\" -> \"int (\" (first occurrence)\n#  2. \", ie. \" -> \", i.e. \" (congestion-window paragraph)\n#  3. \"...the window, ie. if the window is 5 ...\" -> \"...the window, i.e. if the window is 5 ...\"\n#  4. \"1 / int(2.5)\" -> \"1 / int (2.5)\"\n#  5. \"Window size at time 21 is 2.75\" -> \"Window size at time 21 is 3\"\n#  6. Empty \"Ans:\" paragraph (time-23 question) gets the full answer appended\n#  7. \"Window size at time31 is 3.5.\" -> \"Window size at time 31 is 3.99\" (and\n#     the _GoBack bookmark moves from the trailing empty paragraph into this run)\n\n$d = $word.ActiveDocument\n\n# ---- 1 & 4. \"int(\" -> \"int (\" everywhere (only the two expected occurrences exist) ----\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"int(\"\n$find.Replacement.Text = \"int (\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# ---- 2 & 3. \"ie.\" -> \"i.e.\" everywhere (only the two expected occurrences exist) ----\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"ie.\"\n$find.Replacement.Text = \"i.e.\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# ---- 5. \"Window size at time 21 is 2.75\" -> \"Window size at time 21 is 3\" ----\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Window size at time 21 is 2.75\"\n$find.Replacement.Text = \"Window size at time 21 is 3\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# ---- 6. Fill in the empty \"Ans:\" paragraph (the one right after the \"time 23\" question) ----\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -eq \"Ans:\" + [char]13) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -gt 0) {\n    $target = $d.Paragraphs.Item($targetIndex)\n\n    # The paragraph mark stops being bold once non-bold content follows \"Ans: \" -\n    # set the whole paragraph/mark to non-bold first, then re-bold just \"Ans:\".\n    $target.Range.Font.Bold = 0\n\n    $scoped = $d.Range($target.Range.Start, $target.Range.End)\n    $ansFind = $scoped.Find\n    $ansFind.ClearFormatting()\n    $ansFind.Text = \"Ans:\"\n    $ansFind.Execute() | Out-Null\n    $ansFind.Parent.Font.Bold = 1\n\n    $target2 = $d.Paragraphs.Item($targetIndex)\n    $insertPoint = $d.Range($target2.Range.Start, $target2.Range.End - 1)\n    $insertPoint.InsertAfter(\" \")\n    $insertPoint.Font.Bold = 1\n\n    $target3 = $d.Paragraphs.Item($targetIndex)\n    $insertPoint2 = $d.Range($target3.Range.Start, $target3.Range.End - 1)\n    $answerText = \"At time 23, the window size is 3.33 and the number of packets in the network is 2. \" +\n        \"Therefore, one more packet can be pushed into the network as that can be fully accommodated in the congestion window. \" +\n        \"That\" + [char]0x2019 + \"s the reason why a packet is pushed out at time 23 even though no ack is received at that time.\"\n    $insertPoint2.InsertAfter($answerText)\n\n    $target4 = $d.Paragraphs.Item($targetIndex)\n    $scoped2 = $d.Range($target4.Range.Start, $target4.Range.End)\n    $ansFind2 = $scoped2.Find\n    $ansFind2.ClearFormatting()\n    $ansFind2.Text = $answerText\n    $ansFind2.Execute() | Out-Null\n    $ansFind2.Parent.Font.Bold = 0\n}\n\n# ---- 7. \"Window size at time31 is 3.5.\" -> \"Window size at time 31 is 3.99\" + move _GoBack bookmark ----\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Window size at time31 is 3.5.\"\n$find.Replacement.Text = \"Window size at time 31 is 3.99\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# Relocate the _GoBack bookmark from the trailing empty paragraph to just\n# before \"31 is 3.99\".\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"31 is 3.99\"\n$find2.Execute() | Out-Null\n$bmRange = $find2.Parent\n$bmRange.Collapse(1)  # wdCollapseStart\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
